$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Table recompute: shift the chi-square binning table up one "slot"
# (the interval used for the first bin changes from [B18..) based split
# into a finer split anchored on B17/D17, pushing every subsequent bin
# down by one row within the same 28-34 range). ---

$ws.Range("A28").Value = -10000000000
$ws.Range("B28").Formula = '=B17'
$ws.Range("C28").Formula = '=D17'
$ws.Range("D28").Formula = '=NORM.DIST(B28,$H$17,$H$21,TRUE)'
$ws.Range("E28").Formula = '=$O$2*D28'
$ws.Range("F28").Formula = '=C28-$O$2*D28'
$ws.Range("G28").Formula = '=POWER(F28,2)'
$ws.Range("H28").Formula = '=G28/E28'
$ws.Range("I28").Formula = '=(POWER(C28,2))/E28'

$ws.Range("A29").Formula = '=A18'
$ws.Range("B29").Formula = '=B18'
$ws.Range("C29").Formula = '=D18'
$ws.Range("D29").Formula = '=NORM.DIST(B29,$H$17,$H$21,TRUE)-NORM.DIST(A29,$H$17,$H$21,TRUE)'
$ws.Range("E29").Formula = '=$O$2*D29'
$ws.Range("F29").Formula = '=C29-$O$2*D29'
$ws.Range("G29").Formula = '=POWER(F29,2)'
$ws.Range("H29").Formula = '=G29/E29'
$ws.Range("I29").Formula = '=(POWER(C29,2))/E29'

$ws.Range("A30").Formula = '=A19'
$ws.Range("B30").Formula = '=B19'
$ws.Range("C30").Formula = '=D19'
$ws.Range("D30").Formula = '=NORM.DIST(B30,$H$17,$H$21,TRUE)-NORM.DIST(A30,$H$17,$H$21,TRUE)'
$ws.Range("E30").Formula = '=$O$2*D30'
$ws.Range("F30").Formula = '=C30-$O$2*D30'
$ws.Range("G30").Formula = '=POWER(F30,2)'
$ws.Range("H30").Formula = '=G30/E30'
$ws.Range("I30").Formula = '=(POWER(C30,2))/E30'

$ws.Range("A31").Formula = '=A20'
$ws.Range("B31").Formula = '=B20'
$ws.Range("C31").Formula = '=D20'
$ws.Range("D31").Formula = '=NORM.DIST(B31,$H$17,$H$21,TRUE)-NORM.DIST(A31,$H$17,$H$21,TRUE)'
$ws.Range("E31").Formula = '=$O$2*D31'
$ws.Range("F31").Formula = '=C31-$O$2*D31'
$ws.Range("G31").Formula = '=POWER(F31,2)'
$ws.Range("H31").Formula = '=G31/E31'
$ws.Range("I31").Formula = '=(POWER(C31,2))/E31'

$ws.Range("A32").Formula = '=A21'
$ws.Range("B32").Formula = '=B21'
$ws.Range("C32").Formula = '=D21'
$ws.Range("D32").Formula = '=NORM.DIST(B32,$H$17,$H$21,TRUE)-NORM.DIST(A32,$H$17,$H$21,TRUE)'
$ws.Range("E32").Formula = '=$O$2*D32'
$ws.Range("F32").Formula = '=C32-$O$2*D32'
$ws.Range("G32").Formula = '=POWER(F32,2)'
$ws.Range("H32").Formula = '=G32/E32'
$ws.Range("I32").Formula = '=(POWER(C32,2))/E32'

$ws.Range("A33").Formula = '=A22'
$ws.Range("B33").Formula = '=B22'
$ws.Range("C33").Formula = '=D22'
$ws.Range("D33").Formula = '=1-NORM.DIST(A33,$H$17,$H$21,TRUE)'
$ws.Range("E33").Formula = '=$O$2*D33'
$ws.Range("F33").Formula = '=C33-$O$2*D33'
$ws.Range("G33").Formula = '=POWER(F33,2)'
$ws.Range("H33").Formula = '=G33/E33'
$ws.Range("I33").Formula = '=(POWER(C33,2))/E33'

$ws.Range("A34").Formula = '=A23'
$ws.Range("B34").Formula = '=B23'
$ws.Range("C34").Formula = '=D23'
$ws.Range("D34:I34").ClearContents()

# --- Small formatting tweak: the "округляем" label becomes bold ---
$ws.Range("E14").Font.Bold = $true

# --- Restore the active selection to D25 ---
[void]$ws.Range("D25").Select()
